$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 509913.7
Write-Output $ws.Range("H6").Value
